$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "NewsTypes"
$ws.Range("I22").Select()
